# Apply the "eventos_diario" update: append 3 new telemetry event rows
# (rows 2-4) beneath the existing header row on the active sheet.
#
# Columns: A eventId | B Tipo de evento | C Hora | D vehicleId | E Unidad
#          F driverId | G Operador | H latitude | I longitude
#          J maxAcceleration | K video_Interior | L video_Exterior
#
# D/E/F hold digit-only identifiers that must stay textual (as in the
# source export), so their number format is forced to Text ("@") right
# before the value is written; every other column is left alone and
# naturally keeps its text/numeric nature from the assigned value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $val
}

# ---- Row 2 ----
$ws.Range("A2").Value = "281474991205341-1739219899783"
$ws.Range("B2").Value = "No Seat Belt"
$ws.Range("C2").Value = "2025-02-10T14:38:19.783"
Set-TextValue "D2" "281474991205341"
Set-TextValue "E2" "140"
Set-TextValue "F2" "51834149"
$ws.Range("G2").Value = "ABRAHAM ARANA"
$ws.Range("H2").Value = 20.73980557
$ws.Range("I2").Value = -103.39610669
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = "https://s3.samsara.com/samsara-cvdata/4006124/281474991205341/1739219897283/FWyh2Y69B9-camera-video-segment-driver-1739219899783.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSCE7L2VFX%2F20250211%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250211T140017Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjELz%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FwEaCXVzLXdlc3QtMiJHMEUCIQDZTnZFXI0dFVkXIPHvYQML9ifEKfE5s2p0HBJmVerX6wIge5TS7pKkHbsoSAUjnCDmnIoTqvNlErocmN0F%2BZYtjbUq5gMI1f%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FARAEGgw3ODEyMDQ5NDIyNDQiDFzztFRNxsknExqvPiq6A4Hcxe5NYq7HXI73gj4jB%2BYHMYlNgRcfg2mW6Dm7xWNqktoajWlH7N9k3N4jHM9bzc%2BPR5oTgGnN8rbVJjaXTOu7ukvyzT3GRcXF%2F3NnP2Qj%2F1KXF9RgkJVy7bTIvrEkhF8iN1IC7PS9xq260OF3RR2Er2UKIWpPp6EeTTiECwWkDlS4eAiq0dnLsJE7Td33HACfm6y8wH6v1tZ2ZlnMcRQfbU26RdHL62WtGW0zsdJ1TpJeFFIbDh%2F61oIi94L5BSvNoZTcGgsQLsdejeFxxRr0wHJ8P4RPsGkXOXuoyCaVApRLImf6192W4eVfFZ0kbkCYSlcOoNqIFi1V1bwjB2Xu%2B3hASeh9ppLnImzFh0zBa1sL4KhLFG1IhI1EX2NSoqjjhvm22eUI3R0J1JVrbE9TxTu1zr1%2FSK5XlLzGNO6C%2BBsqlSDuEpZrKDCptcggzNVXP18NSuetud7rDkYTX0I%2BqGSDgnDvP8jLKB9%2F%2FKif0N3zf9O5WMlQwqggeZXYm9wqlA91R%2FbQimbe4MfR8VOh6je4yHGo9zTJG5S8vTKahGx9j9YZ9IiNFCKd3pBVBFAXHRwpzlxTVTIwjfCsvQY6pQGcH%2FYS4xeG3optQVxIG%2BEYLnQwXu2fw95Glc4gXbtlUmIRwgagKP0XGbcSkkeuBVqi5NEhUSF8C%2F6XO1Dh4Pi%2FTsWmzfVhZjYypb1Oxb24jCnzOg0WMmbbYEHdIA1UQvtVd%2BsDv7RWlShlwN4lRhvt0aoYFXN9uyXd56%2F%2FqzklYY2l8%2BsURwdc1uV8cdXGvB6ouxvMG%2FCCNrwbEp9iJDZA0N8uU%2BE%3D&X-Amz-SignedHeaders=host&response-expires=Tue%2C%2011%20Feb%202025%2022%3A00%3A17%20GMT&X-Amz-Signature=901e6f43cf8ddb6e3bcabfac97d9ff868e519c066b169e42d21a8088443ca15b"
$ws.Range("L2").Value = "No video URL"

# ---- Row 3 ----
$ws.Range("A3").Value = "281474991205262-1739219492351"
$ws.Range("B3").Value = "No Seat Belt"
$ws.Range("C3").Value = "2025-02-10T14:31:32.351"
Set-TextValue "D3" "281474991205262"
Set-TextValue "E3" "132"
Set-TextValue "F3" "52215867"
$ws.Range("G3").Value = "EMMANUEL SALCEDO"
$ws.Range("H3").Value = 20.57208232
$ws.Range("I3").Value = -103.29531058
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = "https://s3.samsara.com/samsara-cvdata/4006124/281474991205262/1739219489851/3YQOBkkdoW-camera-video-segment-driver-1739219492351.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSCE7L2VFX%2F20250211%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250211T140017Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjELz%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FwEaCXVzLXdlc3QtMiJHMEUCIQDZTnZFXI0dFVkXIPHvYQML9ifEKfE5s2p0HBJmVerX6wIge5TS7pKkHbsoSAUjnCDmnIoTqvNlErocmN0F%2BZYtjbUq5gMI1f%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FARAEGgw3ODEyMDQ5NDIyNDQiDFzztFRNxsknExqvPiq6A4Hcxe5NYq7HXI73gj4jB%2BYHMYlNgRcfg2mW6Dm7xWNqktoajWlH7N9k3N4jHM9bzc%2BPR5oTgGnN8rbVJjaXTOu7ukvyzT3GRcXF%2F3NnP2Qj%2F1KXF9RgkJVy7bTIvrEkhF8iN1IC7PS9xq260OF3RR2Er2UKIWpPp6EeTTiECwWkDlS4eAiq0dnLsJE7Td33HACfm6y8wH6v1tZ2ZlnMcRQfbU26RdHL62WtGW0zsdJ1TpJeFFIbDh%2F61oIi94L5BSvNoZTcGgsQLsdejeFxxRr0wHJ8P4RPsGkXOXuoyCaVApRLImf6192W4eVfFZ0kbkCYSlcOoNqIFi1V1bwjB2Xu%2B3hASeh9ppLnImzFh0zBa1sL4KhLFG1IhI1EX2NSoqjjhvm22eUI3R0J1JVrbE9TxTu1zr1%2FSK5XlLzGNO6C%2BBsqlSDuEpZrKDCptcggzNVXP18NSuetud7rDkYTX0I%2BqGSDgnDvP8jLKB9%2F%2FKif0N3zf9O5WMlQwqggeZXYm9wqlA91R%2FbQimbe4MfR8VOh6je4yHGo9zTJG5S8vTKahGx9j9YZ9IiNFCKd3pBVBFAXHRwpzlxTVTIwjfCsvQY6pQGcH%2FYS4xeG3optQVxIG%2BEYLnQwXu2fw95Glc4gXbtlUmIRwgagKP0XGbcSkkeuBVqi5NEhUSF8C%2F6XO1Dh4Pi%2FTsWmzfVhZjYypb1Oxb24jCnzOg0WMmbbYEHdIA1UQvtVd%2BsDv7RWlShlwN4lRhvt0aoYFXN9uyXd56%2F%2FqzklYY2l8%2BsURwdc1uV8cdXGvB6ouxvMG%2FCCNrwbEp9iJDZA0N8uU%2BE%3D&X-Amz-SignedHeaders=host&response-expires=Tue%2C%2011%20Feb%202025%2022%3A00%3A17%20GMT&X-Amz-Signature=6f027ea24d8cd65b6fd51799edc9ce0122fbb9b298155ebc8abbb5c4851f74cd"
$ws.Range("L3").Value = "No video URL"

# ---- Row 4 ----
$ws.Range("A4").Value = "281474991205262-1739216259756"
$ws.Range("B4").Value = "Harsh Brake"
$ws.Range("C4").Value = "2025-02-10T13:37:39.756"
Set-TextValue "D4" "281474991205262"
Set-TextValue "E4" "132"
Set-TextValue "F4" "52215867"
$ws.Range("G4").Value = "EMMANUEL SALCEDO"
$ws.Range("H4").Value = 20.65361635
$ws.Range("I4").Value = -103.31194159
$ws.Range("J4").Value = 0.7248916029930115
$ws.Range("K4").Value = "https://s3.samsara.com/samsara-dashcam-videos/4006124/281474991205262/1739216254756/Rx1LXQjmB8-camera-video-segment-driver-1739216259756.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSCE7L2VFX%2F20250211%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250211T140017Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjELz%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FwEaCXVzLXdlc3QtMiJHMEUCIQDZTnZFXI0dFVkXIPHvYQML9ifEKfE5s2p0HBJmVerX6wIge5TS7pKkHbsoSAUjnCDmnIoTqvNlErocmN0F%2BZYtjbUq5gMI1f%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FARAEGgw3ODEyMDQ5NDIyNDQiDFzztFRNxsknExqvPiq6A4Hcxe5NYq7HXI73gj4jB%2BYHMYlNgRcfg2mW6Dm7xWNqktoajWlH7N9k3N4jHM9bzc%2BPR5oTgGnN8rbVJjaXTOu7ukvyzT3GRcXF%2F3NnP2Qj%2F1KXF9RgkJVy7bTIvrEkhF8iN1IC7PS9xq260OF3RR2Er2UKIWpPp6EeTTiECwWkDlS4eAiq0dnLsJE7Td33HACfm6y8wH6v1tZ2ZlnMcRQfbU26RdHL62WtGW0zsdJ1TpJeFFIbDh%2F61oIi94L5BSvNoZTcGgsQLsdejeFxxRr0wHJ8P4RPsGkXOXuoyCaVApRLImf6192W4eVfFZ0kbkCYSlcOoNqIFi1V1bwjB2Xu%2B3hASeh9ppLnImzFh0zBa1sL4KhLFG1IhI1EX2NSoqjjhvm22eUI3R0J1JVrbE9TxTu1zr1%2FSK5XlLzGNO6C%2BBsqlSDuEpZrKDCptcggzNVXP18NSuetud7rDkYTX0I%2BqGSDgnDvP8jLKB9%2F%2FKif0N3zf9O5WMlQwqggeZXYm9wqlA91R%2FbQimbe4MfR8VOh6je4yHGo9zTJG5S8vTKahGx9j9YZ9IiNFCKd3pBVBFAXHRwpzlxTVTIwjfCsvQY6pQGcH%2FYS4xeG3optQVxIG%2BEYLnQwXu2fw95Glc4gXbtlUmIRwgagKP0XGbcSkkeuBVqi5NEhUSF8C%2F6XO1Dh4Pi%2FTsWmzfVhZjYypb1Oxb24jCnzOg0WMmbbYEHdIA1UQvtVd%2BsDv7RWlShlwN4lRhvt0aoYFXN9uyXd56%2F%2FqzklYY2l8%2BsURwdc1uV8cdXGvB6ouxvMG%2FCCNrwbEp9iJDZA0N8uU%2BE%3D&X-Amz-SignedHeaders=host&response-expires=Tue%2C%2011%20Feb%202025%2022%3A00%3A17%20GMT&X-Amz-Signature=73b5c94e1d540cf2df8bae6296fda1c888d78e485aeb33dabbfb843d8e1de561"
$ws.Range("L4").Value = "https://s3.samsara.com/samsara-cvdata/4006124/281474991205262/1739216254756/wGRUOWVXLR-camera-video-segment-1739216259756.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSCE7L2VFX%2F20250211%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250211T140017Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjELz%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FwEaCXVzLXdlc3QtMiJHMEUCIQDZTnZFXI0dFVkXIPHvYQML9ifEKfE5s2p0HBJmVerX6wIge5TS7pKkHbsoSAUjnCDmnIoTqvNlErocmN0F%2BZYtjbUq5gMI1f%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FARAEGgw3ODEyMDQ5NDIyNDQiDFzztFRNxsknExqvPiq6A4Hcxe5NYq7HXI73gj4jB%2BYHMYlNgRcfg2mW6Dm7xWNqktoajWlH7N9k3N4jHM9bzc%2BPR5oTgGnN8rbVJjaXTOu7ukvyzT3GRcXF%2F3NnP2Qj%2F1KXF9RgkJVy7bTIvrEkhF8iN1IC7PS9xq260OF3RR2Er2UKIWpPp6EeTTiECwWkDlS4eAiq0dnLsJE7Td33HACfm6y8wH6v1tZ2ZlnMcRQfbU26RdHL62WtGW0zsdJ1TpJeFFIbDh%2F61oIi94L5BSvNoZTcGgsQLsdejeFxxRr0wHJ8P4RPsGkXOXuoyCaVApRLImf6192W4eVfFZ0kbkCYSlcOoNqIFi1V1bwjB2Xu%2B3hASeh9ppLnImzFh0zBa1sL4KhLFG1IhI1EX2NSoqjjhvm22eUI3R0J1JVrbE9TxTu1zr1%2FSK5XlLzGNO6C%2BBsqlSDuEpZrKDCptcggzNVXP18NSuetud7rDkYTX0I%2BqGSDgnDvP8jLKB9%2F%2FKif0N3zf9O5WMlQwqggeZXYm9wqlA91R%2FbQimbe4MfR8VOh6je4yHGo9zTJG5S8vTKahGx9j9YZ9IiNFCKd3pBVBFAXHRwpzlxTVTIwjfCsvQY6pQGcH%2FYS4xeG3optQVxIG%2BEYLnQwXu2fw95Glc4gXbtlUmIRwgagKP0XGbcSkkeuBVqi5NEhUSF8C%2F6XO1Dh4Pi%2FTsWmzfVhZjYypb1Oxb24jCnzOg0WMmbbYEHdIA1UQvtVd%2BsDv7RWlShlwN4lRhvt0aoYFXN9uyXd56%2F%2FqzklYY2l8%2BsURwdc1uV8cdXGvB6ouxvMG%2FCCNrwbEp9iJDZA0N8uU%2BE%3D&X-Amz-SignedHeaders=host&response-expires=Tue%2C%2011%20Feb%202025%2022%3A00%3A17%20GMT&X-Amz-Signature=ca4d861c764c119dec2867874f30dd0b64f13118f05e4f5665833b03254fcb24"
